# Add a new "Italy" test-data sheet, cloned from the existing "Slovakia"
# sheet (same layout/styles), with its own market values, and make it the
# active/selected tab - mirroring how the other country sheets were
# originally authored.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Slovakia")

# Clone "Slovakia" and drop the copy immediately after it.
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "Italy"

# Update the market-specific values on the new sheet.
$newSheet.Range("B2").Value = "Italy Market"
$newSheet.Range("B4").Value = "NGC-3145/T2454/T2453/NGC-3145/T2446"

# Match the selection/active-tab state of the edited workbook: the new
# "Italy" sheet becomes the active tab with B4 selected, while "Slovakia"
# falls back to a full-sheet selection and is no longer the active tab.
$src.Range("A1:XFD1048576").Select()
$newSheet.Activate()
$newSheet.Range("B4").Select()
